$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "28.327.26"
$ws.Cells.Item(2, 5).Value = "  -0.67%  "
$ws.Cells.Item(3, 4).Value = "1.560.91"
$ws.Cells.Item(3, 5).Value = "  -0.17%  "
$ws.Cells.Item(4, 5).Value = "  -0.11%  "
$ws.Cells.Item(5, 4).Value = "210.79"
$ws.Cells.Item(5, 5).Value = "  -0.31%  "
$ws.Cells.Item(6, 5).Value = "  -0.67%  "
$ws.Cells.Item(7, 5).Value = "  -0.14%  "
$ws.Cells.Item(8, 4).Value = "44.38"
$ws.Cells.Item(8, 5).Value = "  -4.31%  "
$ws.Cells.Item(9, 4).Value = "23.56"
$ws.Cells.Item(9, 5).Value = "  -2.19%  "
$ws.Cells.Item(10, 5).Value = "  -1.41%  "
$ws.Cells.Item(11, 5).Value = "  -0.87%  "
$ws.Cells.Item(12, 4).Value = "0.0892"
$ws.Cells.Item(12, 5).Value = "  +0.64%  "
$ws.Cells.Item(13, 4).Value = "1.783.42"
$ws.Cells.Item(13, 5).Value = "  -0.22%  "
$ws.Cells.Item(14, 4).Value = "1.570.99"
$ws.Cells.Item(14, 5).Value = "  +0.41%  "
$ws.Cells.Item(15, 4).Value = "28.322.10"
$ws.Cells.Item(15, 5).Value = "  -0.70%  "
$ws.Cells.Item(16, 5).Value = "  -0.60%  "
$ws.Cells.Item(17, 5).Value = "  -1.43%  "
$ws.Cells.Item(18, 5).Value = "  -1.81%  "
$ws.Cells.Item(19, 4).Value = "227.72"
$ws.Cells.Item(19, 5).Value = "  -0.50%  "
$ws.Cells.Item(20, 4).Value = "7.36"
$ws.Cells.Item(20, 5).Value = "  +0.37%  "
$ws.Cells.Item(21, 4).Value = "0.0₃0678"
$ws.Cells.Item(21, 5).Value = "  -2.07%  "
$ws.Cells.Item(22, 5).Value = "  -0.03%  "
$ws.Cells.Item(23, 4).Value = "3.93"
$ws.Cells.Item(23, 5).Value = "  +1.75%  "
$ws.Cells.Item(24, 4).Value = "8.92"
$ws.Cells.Item(24, 5).Value = "  -2.37%  "
$ws.Cells.Item(25, 5).Value = "  -1.87%  "
$ws.Cells.Item(26, 4).Value = "150.33"
$ws.Cells.Item(26, 5).Value = "  -0.01%  "
$ws.Cells.Item(27, 4).Value = "14.89"
$ws.Cells.Item(27, 5).Value = "  -0.48%  "
$ws.Cells.Item(28, 5).Value = "  -1.60%  "
$ws.Cells.Item(29, 5).Value = "  -0.24%  "
$ws.Cells.Item(30, 5).Value = "  -0.12%  "
$ws.Cells.Item(31, 5).Value = "  +2.27%  "
$ws.Cells.Item(32, 4).Value = "1.06"
$ws.Cells.Item(32, 5).Value = "  -4.19%  "
$ws.Cells.Item(33, 4).Value = "3.18"
$ws.Cells.Item(33, 5).Value = "  -0.80%  "
$ws.Cells.Item(34, 5).Value = "  -0.96%  "
$ws.Cells.Item(35, 4).Value = "1.383.32"
$ws.Cells.Item(35, 5).Value = "  -0.85%  "
$ws.Cells.Item(36, 4).Value = "1.07"
$ws.Cells.Item(36, 5).Value = "  +2.39%  "
$ws.Cells.Item(37, 5).Value = "  -3.20%  "
$ws.Cells.Item(38, 4).Value = "2.34"
$ws.Cells.Item(38, 5).Value = "  -0.57%  "
$ws.Cells.Item(39, 5).Value = "  +1.86%  "
$ws.Cells.Item(40, 5).Value = "  -1.85%  "
$ws.Cells.Item(41, 4).Value = "0.519"
$ws.Cells.Item(41, 5).Value = "  -3.19%  "
$ws.Cells.Item(42, 4).Value = "1.94"
$ws.Cells.Item(42, 5).Value = "  +2.90%  "
$ws.Cells.Item(43, 5).Value = "  -0.11%  "
$ws.Cells.Item(44, 2).Value = "ARBITRUM"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(44, 4).Value = "0.782"
$ws.Cells.Item(44, 5).Value = "  -0.75%  "
$ws.Cells.Item(45, 2).Value = "Kaspa"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(45, 4).Value = "0.0471"
$ws.Cells.Item(45, 5).Value = "  -1.82%  "
$ws.Cells.Item(46, 4).Value = "5.33"
$ws.Cells.Item(46, 5).Value = "  -3.21%  "
$ws.Cells.Item(47, 4).Value = "62.16"
$ws.Cells.Item(47, 5).Value = "  -0.78%  "
$ws.Cells.Item(48, 5).Value = "  -6.32%  "
$ws.Cells.Item(49, 4).Value = "1.696.26"
$ws.Cells.Item(49, 5).Value = "  -0.30%  "
$ws.Cells.Item(50, 4).Value = "85.38"
$ws.Cells.Item(50, 5).Value = "  -0.85%  "
$ws.Cells.Item(51, 5).Value = "  -0.85%  "
